# Update countries & provincias Spain
#
# - Reorders two pairs of countries in the "Pais" table:
#     * "Argelia" now appears before "Nueva Zelanda" / "Libano"
#       (previously it appeared after them).
#     * "Nigeria" now appears before "Georgia" / "Montenegro"
#       (previously it appeared after them).
#   Moving a label up the list shifts the rows below it down by one,
#   so every row from the old position of the moved country down to
#   its new position is rewritten with the data that used to sit one
#   row above it; the moved country's own row gets freshly updated
#   COVID-19 figures.
# - Refreshes the daily case/death counters for several other
#   countries that did not change position (Estados Unidos, Chile,
#   Ecuador, Pakistan, Grecia, Bulgaria).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $pais, $casosTotales, $nuevosCasos, $casosActivos, $recuperados, $casosCriticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 1).Value = $pais
    $ws.Cells.Item($row, 2).Value = $casosTotales
    $ws.Cells.Item($row, 3).Value = $nuevosCasos
    $ws.Cells.Item($row, 4).Value = $casosActivos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $casosCriticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# Simple statistic refreshes (no row movement)
Set-Row 4  "Estados Unidos" 105778 1652 2538  101509 2494 35 1731
Set-Row 29 "Chile"          1909   299 61     1842   7    1  6
Set-Row 30 "Ecuador"        1823   196 3      1772   58   7  48
Set-Row 35 "Pakistan"       1420   47  29     1379   7    1  12
Set-Row 43 "Grecia"         1061   95  52     977    66   4  32
Set-Row 71 "Bulgaria"       313    20  9      298    8    3  6

# "Argelia" moves up, ahead of "Nueva Zelanda" and "Libano"
Set-Row 64 "Argelia"        454 45 29 396 0 3 29
Set-Row 65 "Nueva Zelanda"  451 83 50 401 2 0 0
Set-Row 66 "Libano"         412 21 27 377 3 0 8

# "Nigeria" moves up, ahead of "Georgia" and "Montenegro"
Set-Row 113 "Nigeria"       89 19 3  85 0 0 1
Set-Row 114 "Georgia"       85 2  14 71 1 0 0
Set-Row 115 "Montenegro"    82 0  0  81 1 0 1
